$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A100").Value = "Vị trí địa sinh cung Mệnh tại Sinh địa"
$ws.Range("B100").Value = "Vị trí địa sinh cung Mệnh tại Sinh Địa"

$ws.Range("A101").Value = "Vị trí địa sinh cung Mệnh tại Vượng địa"
$ws.Range("B101").Value = "Vị trí địa sinh cung Mệnh tại Vượng Địa"

$ws.Range("A102").Value = "Vị trí địa sinh cung Mệnh tại Bại địa"
$ws.Range("B102").Value = "Vị trí địa sinh cung Mệnh tại Bại địa"

$ws.Range("A103").Value = "Vị trí địa sinh cung Mệnh tại Tuyệt địa"
$ws.Range("B103").Value = "Vị trí địa sinh cung Mệnh tại Tuyệt địa"

$ws.Range("A104").Value = "Vị trí địa sinh cung Mệnh tại Bình thường"
$ws.Range("B104").Value = "Vị trí địa sinh cung Mệnh tại Bình Thường"

$ws.Range("A100:B104").Interior.Color = $ws.Range("A99").Interior.Color

$ws.Range("A104").Select()
$excel.ActiveWindow.ScrollRow = 73
$excel.ActiveWindow.ScrollColumn = 1
